$wb = $excel.ActiveWorkbook

# --- Existing sheet: ValidLogin ---
$ws1 = $wb.Worksheets.Item("ValidLogin")

# --- Add new sheet "InvalidLogin" right after "ValidLogin" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "InvalidLogin"

# Populate test data for invalid login
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xyz"

# --- Update selections / active sheet ---
# ValidLogin: selection becomes A1:B2, no longer the active/selected tab
$ws1.Range("A1:B2").Select() | Out-Null

# InvalidLogin: becomes active sheet with B2 selected
$ws2.Activate() | Out-Null
$ws2.Range("B2").Select() | Out-Null
